$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# All 8 sheets (UK, Belgium, Denmark, Sweden, Norway, Turkey, Italy, Spain)
# share the same template: rows 13-14 ("Wg" / "Accessories") become four
# rows 13-16 ("MX-BBX" / "MX-DPBX" / "Wg" / "Accessories"), adding two new
# test-case rows. The UK sheet additionally gets a new value in B4.
# ---------------------------------------------------------------------------

# --- UK sheet first, so the new shared strings are appended in the same
#     order as the source edit (NGC-3003/T3834, then MX-BBX, then MX-DPBX).
$ws1 = $wb.Worksheets.Item("UK")

$ws1.Range("A13").Copy()
$ws1.Range("A14:A16").PasteSpecial(-4122)

$ws1.Range("B4").Value = "NGC-3003/T3834"
$ws1.Range("A13").Value = "MX-BBX"
$ws1.Range("A14").Value = "MX-DPBX"
$ws1.Range("A15").Value = "Wg"
$ws1.Range("A16").Value = "Accessories"

# --- Remaining market sheets get the same row 13-16 treatment.
foreach ($name in @("Belgium", "Denmark", "Sweden", "Norway", "Turkey", "Italy", "Spain")) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A13").Copy()
    $ws.Range("A14:A16").PasteSpecial(-4122)

    $ws.Range("A13").Value = "MX-BBX"
    $ws.Range("A14").Value = "MX-DPBX"
    $ws.Range("A15").Value = "Wg"
    $ws.Range("A16").Value = "Accessories"
}

# ---------------------------------------------------------------------------
# Selection / active-sheet state.
# Every non-UK sheet ends up with A7:A16 selected (active cell A7); the UK
# sheet is selected last so it remains the active tab, with B4 selected.
# ---------------------------------------------------------------------------
foreach ($name in @("Belgium", "Denmark", "Sweden", "Norway", "Turkey", "Italy", "Spain")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A7:A16").Select() | Out-Null
}

$ws1.Range("B4").Select() | Out-Null
